$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-12-13"

# Update the December row label
$ws.Range("A13").Value = "December (through 12-13)"

# Update December figures (row 13)
$ws.Range("B13").Value = 13
$ws.Range("C13").Value = 38
$ws.Range("D13").Value = 44
$ws.Range("E13").Value = 28
$ws.Range("F13").Value = 22
$ws.Range("G13").Value = 66
$ws.Range("H13").Value = 96
$ws.Range("I13").Value = 57

# Update Total figures (row 14)
$ws.Range("B14").Value = 304
$ws.Range("C14").Value = 601
$ws.Range("D14").Value = 865
$ws.Range("E14").Value = 710
$ws.Range("F14").Value = 556
$ws.Range("G14").Value = 1330
$ws.Range("H14").Value = 1739
$ws.Range("I14").Value = 1573
